# genlog.xlsx edit: rename "genDbChangeLog" -> "generateChangeLog" and add
# a "genType":"data" field to the file/target JSON sample text.
# Also nudges a couple of row heights and the saved selection, matching
# what Excel records after the author reworked the sample sheet while
# adding GenDbChangeLogCommand / UpdateDbChangeLogCommand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (B2:D2): "genDbChangeLog" -> "generateChangeLog" ---------------
$ws.Range("B2:D2").Value = "generateChangeLog"

# --- Row 3 / column B: append a "genType":"data" member to the JSON -------
$full = '{"type":"file","target":"C:\\projs\\auto-test\\src\\test\\resources\\exportDb\\db.conf","path":"work/download/gen1.xml","genType":"data"}'
$ws.Range("B3").Value = $full

$b3 = $ws.Range("B3")
$b3.Characters(3,4).Font.Name = "Microsoft YaHei"
$b3.Characters(3,4).Font.Size = 10
$b3.Characters(7,1).Font.Name = "Meiryo UI"
$b3.Characters(7,1).Font.Size = 10
$b3.Characters(8,1).Font.Name = "Microsoft YaHei"
$b3.Characters(8,1).Font.Size = 10
$b3.Characters(9,88).Font.Name = "Meiryo UI"
$b3.Characters(9,88).Font.Size = 10
$b3.Characters(97,5).Font.Name = "Microsoft YaHei"
$b3.Characters(97,5).Font.Size = 10
$b3.Characters(102,14).Font.Name = "Meiryo UI"
$b3.Characters(102,14).Font.Size = 10
$b3.Characters(116,3).Font.Name = "Microsoft YaHei"
$b3.Characters(116,3).Font.Size = 10
$b3.Characters(119,19).Font.Name = "Meiryo UI"
$b3.Characters(119,19).Font.Size = 10

# --- Row heights ------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 16.5
$ws.Rows.Item(2).RowHeight = 14.25
$ws.Rows.Item(4).RowHeight = 16.5

# --- Saved selection: D10 -> B4 --------------------------------------------
$ws.Range("B4").Select() | Out-Null
